$wb = $excel.ActiveWorkbook

$tag = $wb.Worksheets.Item("TAG")
$tag.Range("B2").Value = 'TNRResult.addSTEPGRP("ONGLET '
$tag.Range("D2").Value = 'TNRResult.addSTEPGRP("ONGLET '
$tag.Range("F2").Value = 'TNRResult.addSTEPGRP("ONGLET '
$tag.Range("B3").Value = 'TNRResult.addSTEPBLOCK("'
$tag.Range("D3").Value = 'TNRResult.addSTEPBLOCK("'
$tag.Range("F3").Value = 'TNRResult.addSTEPBLOCK("'

$tag.Activate()
